$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: "Deskripsi" (J1) and "Stok" (K1)
$ws.Range("J1").Value = "Deskripsi"
$ws.Range("K1").Value = "Stok"

# Give the new header cells the same look (bold, centered, bordered) as the
# rest of the header row by copying the formatting from an existing header
# cell - this reuses the existing cell style instead of creating new ones.
$null = $ws.Range("A1").Copy()
$null = $ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats

# Extend the bordered body formatting used throughout the data rows down the
# two new columns for rows 2-4.
$null = $ws.Range("A2").Copy()
$null = $ws.Range("J2:K4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Restore the cell selection recorded at save time.
$null = $ws.Range("J6").Select()
